$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings are kept as exact text (matching source formatting)
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.482.80'
$ws.Range('E2').Value = '  +1.91%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.662.71'
$ws.Range('E3').Value = '  +1.12%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.26'
$ws.Range('E5').Value = '  +0.72%  '
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4632'
$ws.Range('E7').Value = '  -3.45%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06119'
$ws.Range('E9').Value = '  -0.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.663.40'
$ws.Range('E10').Value = '  +1.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06946'
$ws.Range('E11').Value = '  -1.90%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.62'
$ws.Range('E12').Value = '  -0.41%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.343'
$ws.Range('E13').Value = '  -0.51%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5741'
$ws.Range('E14').Value = '  -4.27%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '75.08'
$ws.Range('E15').Value = '  +1.84%  '
$ws.Range('E16').Value = '  +0.08%  '
$ws.Range('E17').Value = '  +0.15%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '25.486.65'
$ws.Range('E18').Value = '  +1.95%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000006700'
$ws.Range('E19').Value = '  +1.47%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.38'
$ws.Range('E20').Value = '  +0.87%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.877.62'
$ws.Range('E21').Value = '  +1.47%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.421'
$ws.Range('E22').Value = '  +1.37%  '
$ws.Range('E23').Value = '  +0.79%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.218'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '135.16'
$ws.Range('E25').Value = '  +0.72%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '14.88'
$ws.Range('E26').Value = '  -0.22%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.365'
$ws.Range('E27').Value = '  -1.51%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.708'
$ws.Range('E28').Value = '  +3.60%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '104.02'
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.959'
$ws.Range('E30').Value = '  +1.96%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.07712'
$ws.Range('E31').Value = '  +0.40%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.595'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04340'
$ws.Range('E33').Value = '  +1.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.614'
$ws.Range('E34').Value = '  +1.85%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9407'
$ws.Range('E35').Value = '  +1.21%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.5991'
$ws.Range('E36').Value = '  +1.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9090'
$ws.Range('E37').Value = '  +6.90%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.466'
$ws.Range('E38').Value = '  -4.20%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '107.49'
$ws.Range('E39').Value = '  +9.10%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9993'
$ws.Range('E40').Value = '  +0.05%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.831'
$ws.Range('E41').Value = '  +3.89%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.01460'
$ws.Range('E42').Value = '  -3.72%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.3709'
$ws.Range('E43').Value = '  +0.11%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.972'
$ws.Range('E44').Value = '  +6.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1106'
$ws.Range('E45').Value = '  +0.91%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.05254'
$ws.Range('E46').Value = '  +1.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.116'
$ws.Range('E47').Value = '  +0.33%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '30.13'
$ws.Range('E48').Value = '  +3.71%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.520'
$ws.Range('E49').Value = '  +4.68%  '
$ws.Range('E50').Value = '  +0.31%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.001'
$ws.Range('E51').Value = '  +0.33%  '
